$p = $ppt.ActivePresentation

# Slide 2 -> group "Groupe 8" (shape 2 on the slide) -> "Rectangle : avec
# coins arrondis en haut 4" (2nd item in the group) holds the team member
# name "Alicia LACÔTE".
$s = $p.Slides.Item(2)
$grp = $s.Shapes.Item(2)
$sh = $grp.GroupItems.Item(2)
$tr = $sh.TextFrame.TextRange

# The "Ô" sits at position 14 (1-based) within the whole text-frame range
# ("Alicia LACÔTE" is the 4th paragraph, offset by the 3 preceding empty
# paragraphs). Retype just that single character as a plain "O", which is
# exactly what splits the original single run into three runs around the
# changed character.
$oChar = $tr.Characters(14, 1)
$oChar.Text = "O"
